$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data-driven update: each entry is (row, A, B, C, D, E) for the Hydro Production sheet.
# The commit shifts the whole time series forward by 3 days (switch to Summer time / DST),
# so timestamps in column A move from 2025-03-22 to 2025-03-25, and the B/C/D/E readings
# are replaced with the corresponding new values.
$data = @(
    @(2, 45744.01041666666, 38, 1149, 1454.169889819565, 1187),
    @(3, 45744.02083333334, 38, 1130, 1446.246087580676, 1168),
    @(4, 45744.03125, 38, 1116, 1438.322285341786, 1154),
    @(5, 45744.04166666666, 38, 1114, 1430.398483102897, 1152),
    @(6, 45744.05208333334, 38, 1111, 1433.599489305362, 1149),
    @(7, 45744.0625, 38, 1107, 1436.800495507827, 1145),
    @(8, 45744.07291666666, 38, 1110, 1440.001501710293, 1148),
    @(9, 45744.08333333334, 38, 1035, 1443.202507912758, 1073),
    @(10, 45744.09375, 38, 1029, 1442.32091931939, 1067),
    @(11, 45744.10416666666, 38, 1026, 1441.439330726022, 1064),
    @(12, 45744.11458333334, 37, 1026, 1440.557742132655, 1063),
    @(13, 45744.125, 38, 1027, 1439.676153539287, 1065),
    @(14, 45744.13541666666, 41, 1024, 1450.685515855209, 1065),
    @(15, 45744.14583333334, 41, 1023, 1461.69487817113, 1064),
    @(16, 45744.15625, 41, 976, 1472.704240487051, 1017),
    @(17, 45744.16666666666, 41, 973, 1483.713602802973, 1014),
    @(18, 45744.17708333334, 41, 982, 1499.823584838383, 1023),
    @(19, 45744.1875, 41, 982, 1515.933566873794, 1023),
    @(20, 45744.19791666666, 42, 1028, 1532.043548909204, 1070),
    @(21, 45744.20833333334, 45, 1032, 1548.153530944614, 1077),
    @(22, 45744.21875, 86, 1137, 1603.987475197517, 1223),
    @(23, 45744.22916666666, 112, 1138, 1659.821419450421, 1250),
    @(24, 45744.23958333334, 119, 1140, 1715.655363703324, 1259),
    @(25, 45744.25, 127, 1145, 1771.489307956227, 1272),
    @(26, 45744.26041666666, 235, 1420, 1765.664526177909, 1655),
    @(27, 45744.27083333334, 248, 1527, 1759.839744399592, 1775),
    @(28, 45744.28125, 256, 1556, 1754.014962621274, 1812),
    @(29, 45744.29166666666, 260, 1558, 1748.190180842957, 1818),
    @(30, 45744.30208333334, 306, 1499, 1709.883057436347, 1805),
    @(31, 45744.3125, 297, 1495, 1671.575934029737, 1792),
    @(32, 45744.32291666666, 291, 1558, 1633.268810623127, 1849),
    @(33, 45744.33333333334, 301, 1594, 1594.961687216517, 1895),
    @(34, 45744.34375, 185, 1588, 1562.03855129149, 1773),
    @(35, 45744.35416666666, 281, 1567, 1529.115415366462, 1848),
    @(36, 45744.36458333334, 286, 1566, 1496.192279441435, 1852),
    @(37, 45744.375, 292, 1567, 1463.269143516408, 1859),
    @(38, 45744.38541666666, 390, 1396, 1434.376126875796, 1786),
    @(39, 45744.39583333334, 411, 1379, 1405.483110235185, 1790),
    @(40, 45744.40625, 351, 1385, 1376.590093594573, 1736),
    @(41, 45744.41666666666, 347, 1359, 1347.697076953961, 1706),
    @(42, 45744.42708333334, 326, 1342, 1357.593958425939, 1668),
    @(43, 45744.4375, 328, 1314, 1367.490839897918, 1642),
    @(44, 45744.44791666666, 271, 1303, 1377.387721369897, 1574),
    @(45, 45744.45833333334, 245, 1289, 1387.284602841876, 1534),
    @(46, 45744.46875, 205, 1363, 1393.07789931308, 1568),
    @(47, 45744.47916666666, 196, 1359, 1398.871195784284, 1555),
    @(48, 45744.48958333334, 195, 1358, 1404.664492255488, 1553),
    @(49, 45744.5, 137, 1286, 1410.457788726692, 1423),
    @(50, 45744.51041666666, 126, 1251, 1417.510497474614, 1377),
    @(51, 45744.52083333334, 98, 1099, 1424.563206222536, 1197),
    @(52, 45744.53125, 98, 1087, 1431.615914970458, 1185),
    @(53, 45744.54166666666, 102, 1093, 1438.668623718379, 1195),
    @(54, 45744.55208333334, 110, 1127, 1449.63600562514, 1237),
    @(55, 45744.5625, 110, 1135, 1460.603387531901, 1245),
    @(56, 45744.57291666666, 131, 1182, 1471.570769438662, 1313),
    @(57, 45744.58333333334, 0, 0, 1482.538151345423, 0),
    @(58, 45744.59375, 0, 0, 1516.269410147637, 0),
    @(59, 45744.60416666666, 0, 0, 1550.000668949852, 0),
    @(60, 45744.61458333334, 0, 0, 1583.731927752066, 0),
    @(61, 45744.625, 0, 0, 1617.463186554281, 0),
    @(62, 45744.63541666666, 0, 0, 1668.070569860257, 0),
    @(63, 45744.64583333334, 0, 0, 1718.677953166233, 0),
    @(64, 45744.65625, 0, 0, 1769.28533647221, 0),
    @(65, 45744.66666666666, 0, 0, 1819.892719778186, 0),
    @(66, 45744.67708333334, 0, 0, 1874.855570539866, 0),
    @(67, 45744.6875, 0, 0, 1929.818421301547, 0),
    @(68, 45744.69791666666, 0, 0, 1984.781272063226, 0),
    @(69, 45744.70833333334, 0, 0, 2039.744122824907, 0),
    @(70, 45744.71875, 0, 0, 2100.951559457048, 0),
    @(71, 45744.72916666666, 0, 0, 2162.15899608919, 0),
    @(72, 45744.73958333334, 0, 0, 2223.366432721332, 0),
    @(73, 45744.75, 0, 0, 2284.573869353474, 0),
    @(74, 45744.76041666666, 0, 0, 2281.362368048549, 0),
    @(75, 45744.77083333334, 0, 0, 2278.150866743625, 0),
    @(76, 45744.78125, 0, 0, 2274.939365438701, 0),
    @(77, 45744.79166666666, 0, 0, 2271.727864133777, 0),
    @(78, 45744.80208333334, 0, 0, 2252.248954259103, 0),
    @(79, 45744.8125, 0, 0, 2232.77004438443, 0),
    @(80, 45744.82291666666, 0, 0, 2213.291134509757, 0),
    @(81, 45744.83333333334, 0, 0, 2193.812224635084, 0),
    @(82, 45744.84375, 0, 0, 2148.357936560147, 0),
    @(83, 45744.85416666666, 0, 0, 2102.903648485209, 0),
    @(84, 45744.86458333334, 0, 0, 2057.449360410273, 0),
    @(85, 45744.875, 0, 0, 2011.995072335336, 0),
    @(86, 45744.88541666666, 0, 0, 1943.472549397906, 0),
    @(87, 45744.89583333334, 0, 0, 1874.950026460476, 0),
    @(88, 45744.90625, 0, 0, 1806.427503523047, 0),
    @(89, 45744.91666666666, 0, 0, 1737.904980585617, 0),
    @(90, 45744.92708333334, 0, 0, 1698.055077140555, 0),
    @(91, 45744.9375, 0, 0, 1658.205173695494, 0),
    @(92, 45744.94791666666, 0, 0, 1618.355270250432, 0),
    @(93, 45744.95833333334, 0, 0, 1578.50536680537, 0),
    @(94, 45744.96875, 0, 0, 1575.977963224575, 0),
    @(95, 45744.97916666666, 0, 0, 1573.450559643781, 0),
    @(96, 45744.98958333334, 0, 0, 1570.923156062986, 0),
    @(97, 45745, 0, 0, 1568.395752482191, 0)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]
    $ws.Cells.Item($r, 5).Value = $entry[5]
}

Write-Host "Updated $($data.Count) rows"
